# New Changes 12 7 dec
# Updates order/serial/part identifiers on the "input" sheet and moves the
# active selection to D15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input")

# --- Update the header order/serial numbers (row 2) -----------------------
$ws.Range("A2").Value = "Ord1272018"
$ws.Range("O2").Value = "Ser1272018"

# --- Update the repeated serial-number cells in column A ------------------
# (rows 5, 7, 9, 11 all held the old serial number "Ser2812"/"Ser1122")
$ws.Range("A5").Value = "Ser1272018"
$ws.Range("A7").Value = "Ser1272018"
$ws.Range("A9").Value = "Ser1272018"
$ws.Range("A11").Value = "Ser1272018"

# Re-typing these cells in real Excel reset their formatting to the plain
# bordered style (no forced left alignment) used elsewhere in the sheet
# (e.g. B7); reproduce that by pasting just the formats across from B7,
# which already carries that exact style.
$ws.Range("B7").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the part / hardware identifiers --------------------------------
$ws.Range("C5").Value = "HW117623"
$ws.Range("A13").Value = "PL546711"
$ws.Range("A15").Value = "PL546711"

# --- Move the active selection to D15 --------------------------------------
$ws.Range("D15").Select()
